$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Add the new "APR" worksheet as the last tab (after PaymentHistory)
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$aprSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$aprSheet.Name = "APR"

# ------------------------------------------------------------------
# Header row — write in the same order the strings were first
# introduced in the authored workbook (APR TIERS, Min APR, Max APR,
# Min Score, Max Score) so the shared-string table indices line up.
# ------------------------------------------------------------------
$aprSheet.Range("A1").Value = "APR TIERS"
$aprSheet.Range("D1").Value = "Min APR"
$aprSheet.Range("E1").Value = "Max APR"
$aprSheet.Range("B1").Value = "Min Score"
$aprSheet.Range("C1").Value = "Max Score"

# ------------------------------------------------------------------
# Data rows (HIGH / MEDIUM / LOW tiers)
# ------------------------------------------------------------------
$aprSheet.Range("A2").Value = "HIGH"
$aprSheet.Range("B2").Value = 1
$aprSheet.Range("C2").Value = 49
$aprSheet.Range("D2").Value = 19
$aprSheet.Range("E2").Value = 30

$aprSheet.Range("A3").Value = "MEDIUM"
$aprSheet.Range("B3").Value = 50
$aprSheet.Range("C3").Value = 79
$aprSheet.Range("D3").Value = 13
$aprSheet.Range("E3").Value = 18

$aprSheet.Range("A4").Value = "LOW"
$aprSheet.Range("B4").Value = 80
$aprSheet.Range("C4").Value = 100
$aprSheet.Range("D4").Value = 8
$aprSheet.Range("E4").Value = 12

# ------------------------------------------------------------------
# Formatting — reuse the "ConfidenceScore" sheet's cyan palette
# (fontId/fillId/borderId combination) rather than build it by hand.
# ------------------------------------------------------------------
$confSheet = $wb.Worksheets.Item("ConfidenceScore")

# A:B columns -> identical layout to ConfidenceScore A:B
$confSheet.Range("A1:B4").Copy()
$aprSheet.Range("A1:B4").PasteSpecial(-4122) # xlPasteFormats

# D and E columns mirror ConfidenceScore's "B" (value) column format
$confSheet.Range("B1:B4").Copy()
$aprSheet.Range("D1:D4").PasteSpecial(-4122)
$aprSheet.Range("E1:E4").PasteSpecial(-4122)

# C column: data rows use the same style as column B/D/E; the header
# cell is the same fill/border as B1 but bold (creates the new xf).
$confSheet.Range("B2:B4").Copy()
$aprSheet.Range("C2:C4").PasteSpecial(-4122)
$confSheet.Range("B1").Copy()
$aprSheet.Range("C1").PasteSpecial(-4122)
$aprSheet.Range("C1").Font.Bold = $true

$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# Column widths (closest representable values in this engine)
# ------------------------------------------------------------------
$aprSheet.Range("A1:C1").ColumnWidth = 24.666666666666668
$aprSheet.Range("D1").ColumnWidth = 16.333333333333332
$aprSheet.Range("E1").ColumnWidth = 16.666666666666668

# ------------------------------------------------------------------
# Update selections on the pre-existing sheets
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("PaymentGoal")
$ws1.Activate()
$ws1.Range("B13").Select()

$ws2 = $wb.Worksheets.Item("ConfidenceScore")
$ws2.Activate()
$ws2.Range("A1:B4").Select()

$ws3 = $wb.Worksheets.Item("Income")
$ws3.Activate()
$ws3.Range("C1").Select()

$ws4 = $wb.Worksheets.Item("PaymentHistory")
$ws4.Activate()
$ws4.Range("A2").Select()

# ------------------------------------------------------------------
# Make the new APR sheet the active tab with its own selection
# ------------------------------------------------------------------
$aprSheet.Activate()
$aprSheet.Range("G11").Select()
